$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1412.2609
$ws.Range("I92").Value = 1375.5
$ws.Range("K92").Value = 1375.5
$ws.Range("M92").Value = -127.5
$ws.Range("H138").Value = 2211.6667
$ws.Range("I138").Value = 1231.7307
$ws.Range("J138").Value = 4759.5
$ws.Range("K138").Value = 3695.1921
$ws.Range("L138").Value = 14278.5
$ws.Range("M138").Value = 1444.8079
$ws.Range("N138").Value = -24558.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10564.328
$ws.Range("I32").Value = 8163.8706
$ws.Range("K32").Value = 8163.8706
$ws.Range("M32").Value = -7876.8706
$ws.Range("H45").Value = 4694.933
$ws.Range("I45").Value = 3323.3333
$ws.Range("K45").Value = 3323.3333
$ws.Range("M45").Value = -2946.3333
$ws.Range("H74").Value = 1972.303
$ws.Range("I74").Value = 1733.069
$ws.Range("K74").Value = 1733.069
$ws.Range("M74").Value = -859.069
$ws.Range("H77").Value = 1972.303
$ws.Range("I77").Value = 1733.069
$ws.Range("K77").Value = 8665.344999999999
$ws.Range("M77").Value = -4297.344999999999
$ws.Range("H102").Value = 583.8333
$ws.Range("I102").Value = 583.8333
$ws.Range("K102").Value = 583.8333
$ws.Range("M102").Value = 1038.1667
$ws.Range("H109").Value = 49833.332
$ws.Range("J109").Value = 49833.332
$ws.Range("L109").Value = 49833.332
$ws.Range("N109").Value = -52607.332
$ws.Range("H110").Value = 1613.0344
$ws.Range("I110").Value = 1795.0416
$ws.Range("J110").Value = 739.4
$ws.Range("K110").Value = 1795.0416
$ws.Range("L110").Value = 739.4
$ws.Range("M110").Value = 249.9584
$ws.Range("N110").Value = -4829.4
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954
$ws.Range("H122").Value = 5718.6978
$ws.Range("I122").Value = 5004.4644
$ws.Range("K122").Value = 15013.3932
$ws.Range("M122").Value = -12563.3932
$ws.Range("H132").Value = 55560884
$ws.Range("I132").Value = 100004790
$ws.Range("J132").Value = 5998.375
$ws.Range("K132").Value = 300014370
$ws.Range("L132").Value = 17995.125
$ws.Range("M132").Value = -300011840
$ws.Range("N132").Value = -23055.125
$ws.Range("H133").Value = 68816
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -75059

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4522.467
$ws.Range("I86").Value = 1261.2
$ws.Range("J86").Value = 11045
$ws.Range("K86").Value = 1261.2
$ws.Range("L86").Value = 11045
$ws.Range("M86").Value = -138.2
$ws.Range("N86").Value = -13291
$ws.Range("H89").Value = 4522.467
$ws.Range("I89").Value = 1261.2
$ws.Range("J89").Value = 11045
$ws.Range("K89").Value = 6306
$ws.Range("L89").Value = 55225
$ws.Range("M89").Value = -690
$ws.Range("N89").Value = -66457
$ws.Range("H107").Value = 1101.1428
$ws.Range("I107").Value = 1120
$ws.Range("J107").Value = 1054
$ws.Range("K107").Value = 1120
$ws.Range("L107").Value = 1054
$ws.Range("M107").Value = 800
$ws.Range("N107").Value = -4894
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4923.25
$ws.Range("I31").Value = 1438.3334
$ws.Range("K31").Value = 1438.3334
$ws.Range("M31").Value = -1143.3334
$ws.Range("H34").Value = 4923.25
$ws.Range("I34").Value = 1438.3334
$ws.Range("K34").Value = 1438.3334
$ws.Range("M34").Value = -1236.3334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 426.55554
$ws.Range("J92").Value = 512.5
$ws.Range("L92").Value = 1537.5
$ws.Range("N92").Value = -4033.5
$ws.Range("H103").Value = 4198.9
$ws.Range("J103").Value = 4000
$ws.Range("L103").Value = 12000
$ws.Range("N103").Value = -13758
$ws.Range("H114").Value = 643.4286
$ws.Range("I114").Value = 316.44446
$ws.Range("J114").Value = 888.6667
$ws.Range("K114").Value = 949.33338
$ws.Range("L114").Value = 2666.0001
$ws.Range("M114").Value = 2304.66662
$ws.Range("N114").Value = -9174.000100000001
$ws.Range("H129").Value = 2010.0625
$ws.Range("I129").Value = 1346.625
$ws.Range("J129").Value = 2673.5
$ws.Range("K129").Value = 4039.875
$ws.Range("L129").Value = 8020.5
$ws.Range("M129").Value = 960.125
$ws.Range("N129").Value = -18020.5
$ws.Range("H131").Value = 6544.643
$ws.Range("J131").Value = 9078.571
$ws.Range("L131").Value = 27235.713
$ws.Range("N131").Value = -37315.713

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 98073.586
$ws.Range("I70").Value = 281002
$ws.Range("J70").Value = 6609.375
$ws.Range("K70").Value = 281002
$ws.Range("L70").Value = 6609.375
$ws.Range("M70").Value = -280732
$ws.Range("N70").Value = -7149.375
$ws.Range("H73").Value = 98073.586
$ws.Range("I73").Value = 281002
$ws.Range("J73").Value = 6609.375
$ws.Range("K73").Value = 281002
$ws.Range("L73").Value = 6609.375
$ws.Range("M73").Value = -280066
$ws.Range("N73").Value = -8481.375
$ws.Range("H113").Value = 12704.9
$ws.Range("I113").Value = 8012.5
$ws.Range("J113").Value = 15833.167
$ws.Range("K113").Value = 8012.5
$ws.Range("L113").Value = 15833.167
$ws.Range("M113").Value = -5842.5
$ws.Range("N113").Value = -20173.167

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4545.8
$ws.Range("J7").Value = 5933.778
$ws.Range("L7").Value = 5933.778
$ws.Range("N7").Value = -6157.778
$ws.Range("H55").Value = 1093.6786
$ws.Range("I55").Value = 1440
$ws.Range("K55").Value = 1440
$ws.Range("M55").Value = -1267
$ws.Range("H61").Value = 2334.4
$ws.Range("I61").Value = 1689.0834
$ws.Range("J61").Value = 2930.077
$ws.Range("K61").Value = 1689.0834
$ws.Range("L61").Value = 2930.077
$ws.Range("M61").Value = -1487.0834
$ws.Range("N61").Value = -3334.077
$ws.Range("H93").Value = 403423.72
$ws.Range("I93").Value = 3199.7222
$ws.Range("K93").Value = 3199.7222
$ws.Range("M93").Value = -1951.7222
$ws.Range("H100").Value = 107918.91
$ws.Range("I100").Value = 280527
$ws.Range("K100").Value = 280527
$ws.Range("M100").Value = -279986
$ws.Range("H102").Value = 31832.166
$ws.Range("J102").Value = 31832.166
$ws.Range("L102").Value = 31832.166
$ws.Range("N102").Value = -38322.166
$ws.Range("H110").Value = 39492
$ws.Range("J110").Value = 39492
$ws.Range("L110").Value = 39492
$ws.Range("N110").Value = -47672
$ws.Range("H113").Value = 2334.4
$ws.Range("I113").Value = 1689.0834
$ws.Range("J113").Value = 2930.077
$ws.Range("K113").Value = 1689.0834
$ws.Range("L113").Value = 2930.077
$ws.Range("M113").Value = 480.9166
$ws.Range("N113").Value = -7270.077
$ws.Range("H126").Value = 4545.8
$ws.Range("J126").Value = 5933.778
$ws.Range("L126").Value = 17801.334
$ws.Range("N126").Value = -22741.334

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 19999
$ws.Range("J60").Value = 19999
$ws.Range("L60").Value = 19999
$ws.Range("N60").Value = -21643
$ws.Range("H81").Value = 13166.444
$ws.Range("I81").Value = 15799.714
$ws.Range("J81").Value = 3950
$ws.Range("K81").Value = 31599.428
$ws.Range("L81").Value = 7900
$ws.Range("M81").Value = -30538.428
$ws.Range("N81").Value = -10022
$ws.Range("H84").Value = 13166.444
$ws.Range("I84").Value = 15799.714
$ws.Range("J84").Value = 3950
$ws.Range("K84").Value = 157997.14
$ws.Range("L84").Value = 39500
$ws.Range("M84").Value = -152693.14
$ws.Range("N84").Value = -50108
$ws.Range("H86").Value = 100243.5
$ws.Range("J86").Value = 100243.5
$ws.Range("L86").Value = 100243.5
$ws.Range("N86").Value = -102489.5
$ws.Range("H87").Value = 49999.547
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 49999.547
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 49999.547
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -52495.547
$ws.Range("H89").Value = 100243.5
$ws.Range("J89").Value = 100243.5
$ws.Range("L89").Value = 501217.5
$ws.Range("N89").Value = -512449.5
$ws.Range("H90").Value = 49999.547
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 49999.547
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 149998.641
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -162478.641
$ws.Range("H113").Value = 444
$ws.Range("I113").Value = 456.42856
$ws.Range("J113").Value = 422.25
$ws.Range("K113").Value = 1369.28568
$ws.Range("L113").Value = 1266.75
$ws.Range("M113").Value = 800.71432
$ws.Range("N113").Value = -5606.75
